# Update the "timestamp" column (Z) of the log sheet: this is a re-run of the
# logging notebook, so every previously-logged row gets refreshed with the
# timestamp(s) recorded during the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

# Map of row -> new timestamp string, matching the groups written during the
# new run (batches of rows sharing the same timestamp because they were
# logged together).
$timestamps = @{
    2  = "2025-10-17T07:09:28.079880"
    3  = "2025-10-17T07:09:28.079880"
    4  = "2025-10-17T07:09:28.079880"
    5  = "2025-10-17T07:09:28.079880"
    6  = "2025-10-17T07:09:28.079880"
    7  = "2025-10-17T07:09:28.079880"
    8  = "2025-10-17T07:09:28.079880"
    9  = "2025-10-17T07:09:28.079880"
    10 = "2025-10-17T07:09:28.079880"
    11 = "2025-10-17T07:09:28.079880"
    12 = "2025-10-17T07:09:28.079880"
    13 = "2025-10-17T07:09:28.079880"
    14 = "2025-10-17T07:09:28.079880"
    15 = "2025-10-17T07:09:28.079880"
    16 = "2025-10-17T07:09:28.152550"
    17 = "2025-10-17T07:09:28.153548"
    18 = "2025-10-17T07:09:28.153548"
    19 = "2025-10-17T07:09:28.153548"
    20 = "2025-10-17T07:09:28.153548"
    21 = "2025-10-17T07:09:28.153548"
    22 = "2025-10-17T07:09:28.154548"
    23 = "2025-10-17T07:09:28.154548"
    24 = "2025-10-17T07:09:28.154548"
    25 = "2025-10-17T07:09:28.154548"
    26 = "2025-10-17T07:09:28.230356"
    27 = "2025-10-17T07:09:28.230356"
    28 = "2025-10-17T07:09:28.230356"
    29 = "2025-10-17T07:09:28.230356"
    30 = "2025-10-17T07:09:28.230356"
    31 = "2025-10-17T07:09:28.230356"
    32 = "2025-10-17T07:09:28.230356"
    33 = "2025-10-17T07:09:28.230356"
    34 = "2025-10-17T07:09:28.230356"
    35 = "2025-10-17T07:09:28.230356"
    36 = "2025-10-17T07:09:28.230356"
    37 = "2025-10-17T07:09:28.230356"
    38 = "2025-10-17T07:09:28.230356"
    39 = "2025-10-17T07:09:28.230356"
    40 = "2025-10-17T07:09:28.230356"
    41 = "2025-10-17T07:09:28.230356"
    42 = "2025-10-17T07:09:28.230356"
    43 = "2025-10-17T07:09:28.230356"
    44 = "2025-10-17T07:09:28.230356"
    45 = "2025-10-17T07:09:28.230356"
    46 = "2025-10-17T07:09:28.230356"
    47 = "2025-10-17T07:09:28.230356"
    48 = "2025-10-17T07:09:28.230356"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
